$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: merge a run of text that is currently split across several runs
# into a single run, WITHOUT letting the engine's run-coalescing also fold
# in the immediately preceding "." run (which shares identical rPr but must
# stay a separate run per the target document). We do this by briefly
# flipping Bold off on that single "." character (which forces it to be
# serialized as a run with different formatting, so it can't be merged),
# performing the text replace, and then flipping Bold back on -- which
# restores byte-identical formatting/XML for that run.
# ---------------------------------------------------------------------------
function Merge-RunsProtectingDot($mergedText) {
    $probe = $d.Content
    $found = $probe.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $dotPos = $probe.Start - 1
        $dot = $d.Range($dotPos, $probe.Start)
        if ($dot.Text -eq ".") {
            $dot.Bold = 0
        } else {
            $dot = $null
        }
    }

    $target = $d.Content
    $target.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2)

    if ($dot -ne $null) {
        $probe2 = $d.Content
        $probe2.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        $dotPos2 = $probe2.Start - 1
        $dot2 = $d.Range($dotPos2, $probe2.Start)
        $dot2.Bold = 1
    }
}

# 1. controller.loadData(catalog)  -- merge loadData + ( + catalog + ) -> loadData(catalog)
Merge-RunsProtectingDot("loadData(catalog)")

# 2. model.newCatalog()  -- merge new + Catalog() -> newCatalog()
Merge-RunsProtectingDot("newCatalog()")

# 3. model.addBook(catalog, book)  -- merge addBook + ( + catalog, book + ) -> addBook(catalog, book)
Merge-RunsProtectingDot("addBook(catalog, book)")

# 4. model.addTag(catalog, tag)  -- merge addTa + g( + catalog, tag + ) -> addTag(catalog, tag)
Merge-RunsProtectingDot("addTag(catalog, tag)")

# 5. newList(datastructure, cmpfunction, key, filename, delimiter)
#    -- merge newList + (datastructure...) -> single run (neighbours already
#    differ in formatting, no protection needed)
$t5 = "newList(datastructure, cmpfunction, key, filename, delimiter)"
$d.Content.Find.Execute($t5, $true, $false, $false, $false, $false, $true, 1, $false, $t5, 2)

# 6. "lt" + " " + "para obtener..." + "se hace uso de " -> merge the three
#    Dax-Regular runs following "lt" into a single run.
$t6 = " para obtener el elemento que se quiere de la lista especificada. De igual manera, se hace uso de "
$d.Content.Find.Execute($t6, $true, $false, $false, $false, $false, $true, 1, $false, $t6, 2)

# 7. Add the closing observation paragraph text (last, empty paragraph).
$lastP = $d.Paragraphs.Last
$newRunText = "No se observa ningún cambio específico en el  comportamiento del programa con este cambio, lo cual es lo deseable: los TAD lists son creados de manera que sean fácilmente mutables"
$lastP.Range.InsertBefore($newRunText)
$lastP2 = $d.Paragraphs.Last
$lastP2.Range.LanguageID = "es-419"
